# Update the non-standard concepts mapping sheet:
#  - append 17 new source_code_description -> target_concept_id rows
#  - widen column A
#  - convert the populated range into an Excel Table ("Tabelle1") styled
#    with TableStyleLight1 (replaces the old manual banded-fill look)
#  - update the visible selection / scroll position to match the final file

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new rows (19-35) ---------------------------------------------------
$newRows = @(
  @("NIHSS_followup", 42872750),
  @("NIHSS_stroke", 42872750),
  @("stroke severity (e.g. measured by NIHSS)", 42872750),
  @("mRS_follow_up", 3654822),
  @("mRS_stroke", 3654822),
  @("race", 4216292),
  @("income", 4249447),
  @("low socioeconomic status", 4249447),
  @("functional outcome (e.g. measured by mRS)", 3654822),
  @("functional post stroke outcome", 3654822),
  @("elevated CRP", 37108742),
  @("hs cardiac Troponin", 4010497),
  @("nt-pro-BNP", 4186398),
  @("troponin elevation", 4009409),
  @("Second-hand Smoking", 4184633),
  @("smoking", 4298794),
  @("Smoking", 4298794)
)

$row = 19
foreach ($pair in $newRows) {
  $ws.Cells.Item($row, 1).Value = $pair[0]
  $ws.Cells.Item($row, 2).Value = $pair[1]
  $row = $row + 1
}

# --- widen column A -------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 35.3984375

# --- turn A1:B35 into a proper Excel Table --------------------------------
$tableRange = $ws.Range("A1:B35")
$lo = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$lo.Name = "Tabelle1"
$lo.TableStyle = "TableStyleLight1"

# --- restore sheet view to match the saved state ---------------------------
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("E24").Select()
